$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.320.35"
$ws.Range("E2").Value = "  +5.78%  "
$ws.Range("D3").Value = "2.516.28"
$ws.Range("E3").Value = "  +4.07%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Formula = "'324.51"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").Formula = "'105.32"
$ws.Range("E6").Value = "  +3.39%  "
$ws.Range("D7").Formula = "'0.523"
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("D8").Formula = "'0.999"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Formula = "'0.540"
$ws.Range("E9").Value = "  +2.93%  "
$ws.Range("D10").Formula = "'37.07"
$ws.Range("E10").Value = "  +5.07%  "
$ws.Range("E11").Value = "  +2.67%  "
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Formula = "'18.44"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("E14").Value = "  +4.62%  "
$ws.Range("D15").Value = "2.906.13"
$ws.Range("E15").Value = "  +3.93%  "
$ws.Range("D16").Value = "2.498.74"
$ws.Range("E16").Value = "  +4.05%  "
$ws.Range("D17").Formula = "'0.849"
$ws.Range("E17").Value = "  +2.72%  "
$ws.Range("D18").Value = "47.197.16"
$ws.Range("E18").Value = "  +5.90%  "
$ws.Range("D19").Formula = "'12.79"
$ws.Range("E19").Value = "  +4.57%  "
$ws.Range("D20").Formula = "'6.59"
$ws.Range("E20").Value = "  +3.95%  "
$ws.Range("D21").Value = "0.0₃0943"
$ws.Range("E21").Value = "  +2.93%  "
$ws.Range("D22").Formula = "'70.99"
$ws.Range("E22").Value = "  +3.59%  "
$ws.Range("D23").Formula = "'252.22"
$ws.Range("E23").Value = "  +4.02%  "
$ws.Range("D24").Formula = "'2.38"
$ws.Range("E24").Value = "  +5.15%  "
$ws.Range("E25").Value = "  +3.31%  "
$ws.Range("E26").Value = "  +5.60%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Formula = "'2.31"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Formula = "'10.02"
$ws.Range("E29").Value = "  +5.56%  "
$ws.Range("D30").Formula = "'35.21"
$ws.Range("E30").Value = "  +5.36%  "
$ws.Range("E31").Value = "  +7.21%  "
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("D33").Formula = "'19.85"
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("D35").Formula = "'0.0781"
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Formula = "'1.95"
$ws.Range("E37").Value = "  +4.13%  "
$ws.Range("E38").Value = "  +4.47%  "
$ws.Range("E39").Value = "  +4.57%  "
$ws.Range("D40").Formula = "'123.40"
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("D42").Formula = "'2.23"
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("D43").Formula = "'21.58"
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("D44").Formula = "'0.0298"
$ws.Range("E44").Value = "  +3.53%  "
$ws.Range("D45").Value = "1.980.95"
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("E46").Value = "  +3.77%  "
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("D49").Formula = "'9.13"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("E50").Value = "  +17.32%  "
$ws.Range("D51").Formula = "'79.72"
$ws.Range("E51").Value = "  +5.31%  "
